# Insert a new pair of rows (Primera/Segunda) for date 44588 right after the
# header/data already present, pushing all subsequent rows down by 2 (one
# Primera/Segunda pair). The two rows that fall off the bottom of the
# original range become the new last two rows (548/549).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("408:409").Insert()

# New row 408 - Primera, date 44588
$ws.Range("A408").Value = 8
$ws.Range("B408").Value = "Terminal La Palmera de La Serena"
$ws.Range("C408").Value = "Coquimbo"
$ws.Range("D408").Value = 44588
$ws.Range("E408").Value = 4
$ws.Range("F408").Value = 100112008
$ws.Range("G408").Value = "Coliflor"
$ws.Range("H408").Value = "Sin especificar"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 2000
$ws.Range("K408").Value = 850
$ws.Range("L408").Value = 900
$ws.Range("M408").Value = 875
$ws.Range("N408").Value = "$/unidad"
$ws.Range("O408").Value = "Provincia del Elquí"
$ws.Range("P408").Value = 875
$ws.Range("Q408").Value = 1
$ws.Range("R408").Value = "Hortaliza"

# New row 409 - Segunda, date 44588
$ws.Range("A409").Value = 8
$ws.Range("B409").Value = "Terminal La Palmera de La Serena"
$ws.Range("C409").Value = "Coquimbo"
$ws.Range("D409").Value = 44588
$ws.Range("E409").Value = 4
$ws.Range("F409").Value = 100112008
$ws.Range("G409").Value = "Coliflor"
$ws.Range("H409").Value = "Sin especificar"
$ws.Range("I409").Value = "Segunda"
$ws.Range("J409").Value = 1300
$ws.Range("K409").Value = 750
$ws.Range("L409").Value = 800
$ws.Range("M409").Value = 775
$ws.Range("N409").Value = "$/unidad"
$ws.Range("O409").Value = "Provincia del Elquí"
$ws.Range("P409").Value = 775
$ws.Range("Q409").Value = 1
$ws.Range("R409").Value = "Hortaliza"
